# Update the "想去人数" (want-to-go count) figures that were refreshed by the
# site's automated data generator (gh-pages output regenerated at 456a3b4).
#
# Sheet "展览" (Exhibitions) and sheet "全部类型" (All types) both carry the
# same two data rows, so both need the same updates:
#   F2: 260 -> 261
#   F3: 377 -> 378

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 261
    $ws.Range("F3").Value = 378
}
